$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.402.18'
$ws.Range("E2").Value = '  +0.76%  '

$ws.Range("D3").Value = '3.389.31'
$ws.Range("E3").Value = '  -0.16%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.38'
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.47'
$ws.Range("E6").Value = '  -0.83%  '

$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("E8").Value = '  -0.23%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.70'
$ws.Range("E9").Value = '  +2.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.123'
$ws.Range("E10").Value = '  -0.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.387'
$ws.Range("E11").Value = '  -1.79%  '

$ws.Range("D12").Value = '3.975.82'
$ws.Range("E12").Value = '  +0.05%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.125'
$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.15'
$ws.Range("E14").Value = '  -0.61%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000170'
$ws.Range("E15").Value = '  -0.38%  '

$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '3.373.68'
$ws.Range("E16").Value = '  -0.73%  '

$ws.Range("D17").Value = '61.456.17'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.13'
$ws.Range("E18").Value = '  -0.98%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.65'
$ws.Range("E19").Value = '  -1.94%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.95'
$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '389.10'
$ws.Range("E21").Value = '  +0.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.26'
$ws.Range("E22").Value = '  +1.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.554'
$ws.Range("E23").Value = '  -0.78%  '

$ws.Range("E24").Value = '  +0.18%  '

$ws.Range("E25").Value = '  -3.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.191'
$ws.Range("E26").Value = '  +7.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.31'
$ws.Range("E28").Value = '  -1.18%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.01'
$ws.Range("E29").Value = '  +0.64%  '

$ws.Range("E30").Value = '  +0.48%  '

$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.37'
$ws.Range("E32").Value = '  -3.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.39'
$ws.Range("E33").Value = '  -0.68%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.93'
$ws.Range("E34").Value = '  -0.95%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '167.31'
$ws.Range("E35").Value = '  +0.45%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.04'
$ws.Range("E36").Value = '  +1.03%  '

$ws.Range("D37").Value = '3.425.53'
$ws.Range("E37").Value = '  +0.00%  '

$ws.Range("E38").Value = '  -0.44%  '

$ws.Range("E39").Value = '  -0.93%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.98'
$ws.Range("E40").Value = '  -8.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.778'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.43'
$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("E43").Value = '  -0.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.14'
$ws.Range("E44").Value = '  +0.41%  '

$ws.Range("D45").Value = '2.455.23'
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.94'
$ws.Range("E46").Value = '  -1.17%  '

$ws.Range("E47").Value = '  -1.90%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  +0.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0261'
$ws.Range("E49").Value = '  -2.47%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.06'
$ws.Range("E50").Value = '  -1.18%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.207'
$ws.Range("E51").Value = '  -1.38%  '
